# Daily attendance processing - 2025-11-30 11:45:52
# Normalises the "Recorded By" list (column G) on each attendance row:
# the comma-separated recorder names/emails are alphabetised
# (case-insensitive), except that an entry that is the EXACT, literal
# word "System" (capital S) is always moved to the end of the list -
# other casings such as "system" are left in their alphabetical slot.

function Test-ExactMatch {
    param($s, $target)
    if ($s.Length -ne $target.Length) { return $false }
    $sChars = $s.ToCharArray()
    $tChars = $target.ToCharArray()
    for ($i = 0; $i -lt $sChars.Length; $i++) {
        if ([int]$sChars[$i] -ne [int]$tChars[$i]) {
            return $false
        }
    }
    return $true
}

function Get-SortedRecordedBy {
    param($value)

    $parts = $value.Split(",")
    $trimmedParts = @()
    foreach ($p in $parts) {
        $trimmedParts += $p.Trim()
    }

    $nonSystem = @()
    $systemItems = @()
    foreach ($p in $trimmedParts) {
        if (Test-ExactMatch $p "System") {
            $systemItems += $p
        } else {
            $nonSystem += $p
        }
    }

    $nonSystemSorted = $nonSystem | Sort-Object { $_.ToLower() }
    $result = @($nonSystemSorted) + @($systemItems)
    return ($result -join ", ")
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$col = 7  # Column G = "Recorded By"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $current = $cell.Value2
    if ($current -eq $null) { continue }
    if ($current -eq "") { continue }
    if ($current -eq "Recorded By") { continue }

    $updated = Get-SortedRecordedBy $current
    if ($updated -ne $current) {
        $cell.Value2 = $updated
    }
}
